# Apply updated values to Sheet1 for rows 2-25 (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.488102091309003
$ws.Range("C2").Value = 0.04297843819446712
$ws.Range("D2").Value = 0.1370954164261189
$ws.Range("E2").Value = 0.07377937550792169
$ws.Range("F2").Value = 2.443551437920121
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("K2").Value = 1.006837424343303
$ws.Range("L2").Value = 0.2332125897266764
$ws.Range("N2").Value = 3.257016017890948
$ws.Range("B3").Value = 1.437926680809312
$ws.Range("C3").Value = 0.03891360204065109
$ws.Range("D3").Value = 0.1375072156330379
$ws.Range("E3").Value = 0.07328404718238524
$ws.Range("F3").Value = 2.417129626971743
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("K3").Value = 0.9578911728305286
$ws.Range("L3").Value = 0.2262685050182967
$ws.Range("N3").Value = 3.258522478605585
$ws.Range("B4").Value = 1.408023314977669
$ws.Range("C4").Value = 0.03639846191308038
$ws.Range("D4").Value = 0.1377667964015759
$ws.Range("E4").Value = 0.07301236143787904
$ws.Range("F4").Value = 2.401974831123155
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("K4").Value = 0.9284730697622763
$ws.Range("L4").Value = 0.2221438623067087
$ws.Range("N4").Value = 3.260171859267004
$ws.Range("B5").Value = 1.39606465224324
$ws.Range("C5").Value = 0.03536859521879876
$ws.Range("D5").Value = 0.1378742713094931
$ws.Range("E5").Value = 0.07290980921765211
$ws.Range("F5").Value = 2.396067340362109
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("K5").Value = 0.9166442831923405
$ws.Range("L5").Value = 0.220497943314399
$ws.Range("N5").Value = 3.2610258723372
$ws.Range("B6").Value = 1.394092640746152
$ws.Range("C6").Value = 0.03519728704408465
$ws.Range("D6").Value = 0.1378922197359831
$ws.Range("E6").Value = 0.072893273543718
$ws.Range("F6").Value = 2.395102593200065
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("K6").Value = 0.9146897347117147
$ws.Range("L6").Value = 0.2202267462747045
$ws.Range("N6").Value = 3.261178657333545
$ws.Range("B7").Value = 1.4078611166683
$ws.Range("C7").Value = 0.03638459280985273
$ws.Range("D7").Value = 0.1377682389890209
$ws.Range("E7").Value = 0.07301094533160146
$ws.Range("F7").Value = 2.401894075227403
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("K7").Value = 0.9283128979758146
$ws.Range("L7").Value = 0.2221215236095304
$ws.Range("N7").Value = 3.260182640766786
$ws.Range("B8").Value = 1.47061369845494
$ws.Range("C8").Value = 0.04158084816332064
$ws.Range("D8").Value = 0.1372360072341827
$ws.Range("E8").Value = 0.07360185295806865
$ws.Range("F8").Value = 2.43421920214054
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("K8").Value = 0.9898287004302517
$ws.Range("L8").Value = 0.2307893697286545
$ws.Range("N8").Value = 3.257384921364093
$ws.Range("B9").Value = 1.600868671506532
$ws.Range("C9").Value = 0.05162101994606871
$ws.Range("D9").Value = 0.1362457020019416
$ws.Range("E9").Value = 0.07501811626369914
$ws.Range("F9").Value = 2.506110257756134
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("K9").Value = 1.115524746964979
$ws.Range("L9").Value = 0.2488942489683268
$ws.Range("N9").Value = 3.257660754437836
$ws.Range("B10").Value = 1.700995526296367
$ws.Range("C10").Value = 0.05891186340382149
$ws.Range("D10").Value = 0.1355505948547844
$ws.Range("E10").Value = 0.07621589625808056
$ws.Range("F10").Value = 2.564154861022203
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("K10").Value = 1.211005997462451
$ws.Range("L10").Value = 0.2628781551619568
$ws.Range("N10").Value = 3.261398319846379
$ws.Range("B11").Value = 1.74751718180022
$ws.Range("C11").Value = 0.062211293796679
$ws.Range("D11").Value = 0.1352414066785457
$ws.Range("E11").Value = 0.07679502049788312
$ws.Range("F11").Value = 2.59170572017652
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("K11").Value = 1.255133650961881
$ws.Range("L11").Value = 0.2693896310849624
$ws.Range("N11").Value = 3.263871213444816
$ws.Range("B12").Value = 1.765274246038473
$ws.Range("C12").Value = 0.06345831405954527
$ws.Range("D12").Value = 0.1351253347956121
$ws.Range("E12").Value = 0.07701924702134733
$ws.Range("F12").Value = 2.602303973590011
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("K12").Value = 1.271943897005343
$ws.Range("L12").Value = 0.2718770480895927
$ws.Range("N12").Value = 3.264919118760218
$ws.Range("B13").Value = 1.761443698060532
$ws.Range("C13").Value = 0.06318985162285173
$ws.Range("D13").Value = 0.1351502879620021
$ws.Range("E13").Value = 0.07697073679749167
$ws.Range("F13").Value = 2.60001408785871
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("K13").Value = 1.268319052308328
$ws.Range("L13").Value = 0.2713403736907622
$ws.Range("N13").Value = 3.26468846984713
$ws.Range("B14").Value = 1.748975251877766
$ws.Range("C14").Value = 0.06231393463583856
$ws.Range("D14").Value = 0.1352318371043584
$ws.Range("E14").Value = 0.07681336904941816
$ws.Range("F14").Value = 2.592574328078456
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("K14").Value = 1.256514632061851
$ws.Range("L14").Value = 0.2695938378008123
$ws.Range("N14").Value = 3.263955189016514
$ws.Range("B15").Value = 1.741356253126185
$ws.Range("C15").Value = 0.0617770994835638
$ws.Range("D15").Value = 0.1352819200158617
$ws.Range("E15").Value = 0.0767176181079563
$ws.Range("F15").Value = 2.588038808750071
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("K15").Value = 1.249297128114307
$ws.Range("L15").Value = 0.2685268574543898
$ws.Range("N15").Value = 3.26352056154721
$ws.Range("B16").Value = 1.697974843951272
$ws.Range("C16").Value = 0.05869589772535733
$ws.Range("D16").Value = 0.1355709424092595
$ws.Range("E16").Value = 0.07617873827324573
$ws.Range("F16").Value = 2.562377455235023
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("K16").Value = 1.20813612924772
$ws.Range("L16").Value = 0.2624556421503428
$ws.Range("N16").Value = 3.261252289515866
$ws.Range("B17").Value = 1.671611317661757
$ws.Range("C17").Value = 0.05680132354456191
$ws.Range("D17").Value = 0.1357500472107596
$ws.Range("E17").Value = 0.07585692457055515
$ws.Range("F17").Value = 2.546928962713011
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("K17").Value = 1.183062933391881
$ws.Range("L17").Value = 0.2587696516751095
$ws.Range("N17").Value = 3.260058933434053
$ws.Range("B18").Value = 1.656539318896932
$ws.Range("C18").Value = 0.05570999025169954
$ws.Range("D18").Value = 0.1358537241620752
$ws.Range("E18").Value = 0.0756750494266214
$ws.Range("F18").Value = 2.53815124007582
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("K18").Value = 1.16870667239138
$ws.Range("L18").Value = 0.2566636935535485
$ws.Range("N18").Value = 3.259445250397789
$ws.Range("B19").Value = 1.651451921958653
$ws.Range("C19").Value = 0.05534020237146819
$ws.Range("D19").Value = 0.1358889408133521
$ws.Range("E19").Value = 0.07561402328849809
$ws.Range("F19").Value = 2.535197759720859
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("K19").Value = 1.163857070117757
$ws.Range("L19").Value = 0.2559530762531494
$ws.Range("N19").Value = 3.259249943702443
$ws.Range("B20").Value = 1.674408276914733
$ws.Range("C20").Value = 0.05700317178778391
$ws.Range("D20").Value = 0.1357309128192554
$ws.Range("E20").Value = 0.0758908485610128
$ws.Range("F20").Value = 2.548562316091193
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("K20").Value = 1.185725269026449
$ws.Range("L20").Value = 0.2591605691203824
$ws.Range("N20").Value = 3.26017844105985
$ws.Range("B21").Value = 1.752633725649162
$ws.Range("C21").Value = 0.06257127730836487
$ws.Range("D21").Value = 0.135207856693528
$ws.Range("E21").Value = 0.07685945811514827
$ws.Range("F21").Value = 2.594755074238748
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("K21").Value = 1.259979157204668
$ws.Range("L21").Value = 0.2701062492166386
$ws.Range("N21").Value = 3.264167543040458
$ws.Range("B22").Value = 1.804576655317305
$ws.Range("C22").Value = 0.06619638701393171
$ws.Range("D22").Value = 0.1348719028516978
$ws.Range("E22").Value = 0.07752120574669874
$ws.Range("F22").Value = 2.625908701538748
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("K22").Value = 1.309091813167811
$ws.Range("L22").Value = 0.2773861845586509
$ws.Range("N22").Value = 3.267424582822883
$ws.Range("B23").Value = 1.776778765894619
$ws.Range("C23").Value = 0.06426285183080438
$ws.Range("D23").Value = 0.1350506679083132
$ws.Range("E23").Value = 0.07716539201416239
$ws.Range("F23").Value = 2.609193039543896
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("K23").Value = 1.282825929473262
$ws.Range("L23").Value = 0.2734891645519184
$ws.Range("N23").Value = 3.265626647997294
$ws.Range("B24").Value = 1.673143507699137
$ws.Range("C24").Value = 0.0569119227985766
$ws.Range("D24").Value = 0.1357395612742316
$ws.Range("E24").Value = 0.07587550174511648
$ws.Range("F24").Value = 2.547823553671165
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("K24").Value = 1.184521444425684
$ws.Range("L24").Value = 0.2589837943377091
$ws.Range("N24").Value = 3.26012418617772
$ws.Range("B25").Value = 1.564856545663474
$ws.Range("C25").Value = 0.04892033638206783
$ws.Range("D25").Value = 0.1365078984176602
$ws.Range("E25").Value = 0.07460738526533106
$ws.Range("F25").Value = 2.485747211794603
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("K25").Value = 1.080973873212628
$ws.Range("L25").Value = 0.2438771174267771
$ws.Range("N25").Value = 3.256966985954463
